$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "19+39=",
    "53+14=",
    "48-3=",
    "31+27=",
    "39+21=",
    "49+27=",
    "66-38=",
    "39-20=",
    "14+17=",
    "11+38=",
    "51-14=",
    "76-17=",
    "28+70=",
    "55-9=",
    "40+12=",
    "81+7=",
    "63+5=",
    "96-41=",
    "43+53=",
    "47-12=",
    "9+33=",
    "39+12=",
    "40+8=",
    "86-15=",
    "52+45=",
    "99-68=",
    "20+73=",
    "22+20=",
    "74-72=",
    "81-46=",
    "32+67=",
    "16+21=",
    "42-39=",
    "73+8=",
    "23+30=",
    "25+34=",
    "71-57=",
    "53-39=",
    "76-60=",
    "9+51=",
    "37+50=",
    "20-4=",
    "41-41=",
    "25+70=",
    "33-10=",
    "47+8=",
    "97-79=",
    "64-14=",
    "14+37=",
    "22-12=",
    "54-35=",
    "68+15=",
    "75-52=",
    "57+6=",
    "74-2=",
    "84+14=",
    "91-71=",
    "63+19=",
    "58-51=",
    "73-30=",
    "73-46=",
    "41+42=",
    "55-6=",
    "27+58=",
    "82-24=",
    "31+24=",
    "22+5=",
    "24-20=",
    "57-41=",
    "97-90=",
    "71+4=",
    "9+39=",
    "40+25=",
    "32+44=",
    "18+76=",
    "91-73=",
    "21+19=",
    "34+14=",
    "87-51=",
    "29+45=",
    "80-2=",
    "44+0=",
    "20+39=",
    "29-3=",
    "68-5=",
    "24-13=",
    "88-12=",
    "87-37=",
    "38+20=",
    "43-34=",
    "65-37=",
    "77-23=",
    "98-20=",
    "67-66=",
    "4+8=",
    "61-60=",
    "61-37=",
    "62-31=",
    "38+11=",
    "15+7="
)

$cols = 5
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = [Math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    [void]$r.MoveEnd(1, -1)
    $r.Text = $values[$i]
}
Write-Host "Done updating $($values.Length) cells"
